$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
}

Replace-Text "高度で回復力のあるセキュリティ製品です。" "高度で回復性のあるセキュリティ製品です。"
Replace-Text "Firewall Protection:" "ファイアウォールによる保護:"
Replace-Text "アプリケーション層のネットワーク パケットを検査" "アプリケーション レイヤーのネットワーク パケットを検査"
Replace-Text "基づいて規則セットを動的に調整し、アプリケーション層攻撃に関連する" "基づいてルール セットを動的に調整し、アプリケーション レイヤー攻撃に関連する"
Replace-Text "Endpoint Security:" "エンドポイントのセキュリティ:"
Replace-Text "SIEM (セキュリティ情報およびイベント管理)" "SIEM (セキュリティ情報イベント管理)"
Replace-Text "生体認証やスマート カード統合など" "生体認証やスマート カードの統合など"
Replace-Text "SDN ハードウェア要件" "2.1 ハードウェア要件"
Replace-Text " 2.5 GHz 以上のクアッド コア(ハードウェア アクセラレーションサポート付き)" " クアッド コア 2.5 GHz 以上 (ハードウェア アクセラレーション サポート付き)"
Replace-Text " 16 GB 以上、ECC (エラー修正コード) をお勧めします" " 16 GB 以上、ECC (エラー訂正コード) をお勧めします"
Replace-Text "Storage:" "ストレージ:"
Replace-Text "ネットワーク インターフェイス カード (NIC):ジャンボ フレームをサポート" "ネットワーク インターフェイス カード (NIC):"
Replace-Text " デュアル ギガビット イーサネット" " ジャンボ フレームをサポートするデュアル ギガビット イーサネット"
Replace-Text "ソフトウェア要件。" "2.2 ソフトウェア要件"
Replace-Text " Windows Server 2019 以降、CentOS 8 以降と互換性があります" " Windows Server 2019 以降、CentOS 8 以降または同等のものと互換性があります"
Replace-Text "Database:" "データベース:"
Replace-Text " PostgreSQL 13 for Data Storage(高パフォーマンスのインデックス作成用に最適化)" " データ ストレージのための PostgreSQL 13、高パフォーマンスのインデックス作成用に最適化"
Replace-Text "展開の手順" "3.1 デプロイメントの手順"
Replace-Text "展開前評価:" "デプロイ前評価:"
Replace-Text "Installation:" "インストール:"
Replace-Text "メンテナンスとサポート" "3.2 メンテナンスとサポート"

"Done"
